# fixed sorting of plot selector
# The "order" column (T) on the element_info sheet drove the sort order used
# by the plot-selector UI. The previous values were reversed relative to the
# intended sort, so every rank N is replaced by (71 - N) -- i.e. the ranking
# is flipped end-for-end (1<->70, 2<->69, ... 35<->36).
# Also update the view so the previously-selected column/scroll position
# reflects the new "order" column (V) used for sorting, instead of the old
# stale D-column selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("element_info")

# --- Reverse the rank values stored in column T ("order"), rows 2-85 ---
# Column T is the 20th column (A=1 ... T=20). Only rows that actually carry
# an order value are touched; blank cells are left untouched.
$orderCol = 20
for ($r = 2; $r -le 85; $r++) {
    $cell = $ws.Cells.Item($r, $orderCol)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = 71 - $current
    }
}

# --- Update the sheet view: selection now on V2:V85, scrolled to E72 ---
[void]$ws.Activate()

$win = $excel.ActiveWindow
# Scroll the frozen (bottom-right) pane so E72 becomes the top-left visible
# cell, then select V2:V85 (active cell V2) to match the new selector range.
$win.ScrollRow = 72
$win.ScrollColumn = 5
[void]$ws.Range("V2:V85").Select()
